$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    2  = "276.94"
    3  = "21.17"
    4  = "6.268"
    5  = "0.06220"
    6  = "3.556"
    7  = "1.534"
    8  = "6.582"
    9  = "0.8276"
    10 = "0.1665"
    11 = "0.08305"
    12 = "0.03502"
    13 = "0.03166"
    14 = "0.09172"
    15 = "3.761"
    16 = "0.001639"
    17 = "0.04683"
    18 = "0.006397"
    19 = "0.006218"
    21 = "0.0001498"
    23 = "2.312"
    40 = "0.04746"
    41 = "0.007070"
    42 = "0.1121"
    43 = "0.005194"
    44 = "0.01135"
    45 = "0.00006294"
    46 = "0.0009897"
    47 = "0.00000000749"
    48 = "0.9902"
    49 = "0.001400"
    50 = "0.00001898"
    51 = "0.01238"
}

foreach ($row in $changes.Keys) {
    $ws.Range("D$row").Value = "'" + $changes[$row]
}
